# Update the "Corr/total marks" marksheet figures.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("quiz")

# Total correct count (row 11 = "Total" correct marks, was 3 -> 5)
$ws.Range("B11").Value = 5

# Total marks (row 12 "Total" row, B column, was 39 -> 65)
$ws.Range("B12").Value = 65

# Correct/total summary text (E12, was "34/84" -> "65/140")
$ws.Range("E12").Value = "65/140"
